# Insert a new "Page 3" heading paragraph before the first paragraph.
$d = $word.ActiveDocument

$firstParaRange = $d.Paragraphs.Item(1).Range
$firstParaRange.InsertParagraphBefore()

# The inserted paragraph is now paragraph 1; style and fill it in.
$newPara = $d.Paragraphs.Item(1)
$newPara.Style = "Heading1"
$newPara.Range.Text = "Page 3"
